$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the old C1 formula (SUM(B1:B6)) down to C2
$ws.Range("C2").Formula = "=SUM(B1:B6)"

# Move the old B8 formula (SUM(A1:A7)+C1) up to C1, but now it needs to reference
# the relocated cell C2 instead of the old C1
$ws.Range("C1").Formula = "=SUM(A1:A7)+C2"

# A7 used to reference C1 (the SUM(B1:B6) cell); update it to reference the
# relocated C2 cell so the result stays the same
$ws.Range("A7").Formula = "=C2+B5"

# Remove the now-obsolete row 8 (its formula moved into C1)
$ws.Range("B8").ClearContents()

# Update the "Name" defined name to point at the new location of the
# SUM(A1:A7)+C2 formula (now in C1, previously in B8). Delete + re-add so the
# calculation engine picks up the new target.
$wb.Names.Item("Name").Delete()
$wb.Names.Add("Name", "=Sheet1!`$C`$1")

# Update the selected cell to reflect the new active selection
$ws.Range("C1").Select()

$excel.CalculateFullRebuild()

$wb.Save()
